# Updated cryptos list on Sat May 13 15:06:32 UTC 2023 with GitHub Actions
# Refreshes the Price (D) / Volume(1h) (E) columns for existing coins, and
# shifts rows 45-51 so the oldest entry (PaxosStandard) drops off the list
# and a new coin (RenderToken) is appended at the bottom.
#
# For column D, values are forced to Text ("@") before assignment because many
# of the price strings (e.g. "318.20", "157.00") are valid numeric literals
# that Excel would otherwise silently reinterpret as numbers (losing the
# trailing zero). ClearFormats() afterwards removes the now-unneeded number
# format so the cell keeps the workbook's original (default) style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.419.52'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +3.65%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.837.13'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +3.48%  '
$ws.Range("E4").Value = '  +2.39%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.20'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +3.63%  '
$ws.Range("E6").Value = '  +2.20%  '
$ws.Range("E7").Value = '  +1.74%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3724'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +2.56%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07340'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +2.48%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8726'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +3.62%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.39'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +4.47%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.968.35'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +10.00%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.470'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +4.13%  '
$ws.Range("E14").Value = '  +3.57%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07149'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +3.70%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '82.25'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +4.31%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.031'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +2.38%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008981'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +3.19%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.025'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +2.30%  '
$ws.Range("E20").Value = '  +3.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.436.98'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +3.69%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.247'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +2.70%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.18'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.47%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.164.20'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +7.69%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.00'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +2.86%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.892'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.29%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.54'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +2.80%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.243'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +3.37%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.925'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +7.72%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '115.35'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +1.24%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09045'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +1.57%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.197'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +6.76%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7586'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +4.19%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.463'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +3.10%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.861'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +4.22%  '
$ws.Range("E36").Value = '  +2.43%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.151'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +4.00%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01955'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +3.54%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05248'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +2.12%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5166'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +4.84%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.789'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +7.28%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1661'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +2.97%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.514'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +2.71%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.472'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +6.07%  '
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '108.68'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +3.50%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.49'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +3.53%  '
$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.028'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.44%  '
$ws.Range("B48").Value = 'Decentraland'
$ws.Range("C48").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4627'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +3.72%  '
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.668'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +2.21%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06293'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.67%  '
$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.869'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +8.94%  '
